$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# M3: plain text date stamp update
$ws.Range("M3").Value = "Printed On: 10/23/2025"

# B16: Agency name -- preserve the bold "Agency:" label run, replace only the city/borough name
$b16 = $ws.Range("B16")
$b16.Value = "Agency: MCSHERRYSTOWN BOROUGH"
$prefix = $b16.Characters(1, 7)
$prefix.Font.Bold = $true
$prefix.Font.Name = "verdana"
$prefix.Font.Size = 9
$suffix = $b16.Characters(8, $b16.Characters().Text.Length - 7)
$suffix.Font.Bold = $false
$suffix.Font.Name = "verdana"
$suffix.Font.Size = 9

# Data grid cells: all are stored as text (not numbers). Force text number format
# so the value round-trips as a string, and restore the Arial 9 font with the correct
# bold/non-bold weight matching each row (bold on subtotal/total rows).
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("I21")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("J21")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("I22")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("I23")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("J23")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("I24")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("I25")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "8"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("I29")
$c.NumberFormat = "@"
$c.Value = "8"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("J29")
$c.NumberFormat = "@"
$c.Value = "7"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = "4"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("I30")
$c.NumberFormat = "@"
$c.Value = "4"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("J30")
$c.NumberFormat = "@"
$c.Value = "4"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("I31")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("J31")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "4"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("I32")
$c.NumberFormat = "@"
$c.Value = "4"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("J32")
$c.NumberFormat = "@"
$c.Value = "3"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("I33")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("J33")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = "4"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("I34")
$c.NumberFormat = "@"
$c.Value = "4"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("J34")
$c.NumberFormat = "@"
$c.Value = "7"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("J35")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "4"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("I36")
$c.NumberFormat = "@"
$c.Value = "4"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("J36")
$c.NumberFormat = "@"
$c.Value = "6"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "21"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("I38")
$c.NumberFormat = "@"
$c.Value = "21"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("J38")
$c.NumberFormat = "@"
$c.Value = "25"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("N38")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("I43")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("J43")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "33"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("I45")
$c.NumberFormat = "@"
$c.Value = "33"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("J45")
$c.NumberFormat = "@"
$c.Value = "39"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("N45")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "17"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("I46")
$c.NumberFormat = "@"
$c.Value = "17"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("J46")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("N46")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("I47")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("J47")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "6"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("I48")
$c.NumberFormat = "@"
$c.Value = "6"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("J48")
$c.NumberFormat = "@"
$c.Value = "5"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("N48")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("I49")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("J49")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("I50")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("J50")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "5"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("I51")
$c.NumberFormat = "@"
$c.Value = "5"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("J51")
$c.NumberFormat = "@"
$c.Value = "6"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("E52")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("I52")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("J52")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("N52")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("E54")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("I54")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("J54")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("N54")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("E55")
$c.NumberFormat = "@"
$c.Value = "3"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("I55")
$c.NumberFormat = "@"
$c.Value = "3"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("J55")
$c.NumberFormat = "@"
$c.Value = "5"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("N55")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("E61")
$c.NumberFormat = "@"
$c.Value = "3"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("I61")
$c.NumberFormat = "@"
$c.Value = "3"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("J61")
$c.NumberFormat = "@"
$c.Value = "5"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("N61")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("E62")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("I62")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("J62")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("E63")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("I63")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("J63")
$c.NumberFormat = "@"
$c.Value = "3"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("N63")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("E65")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("I65")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("J65")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("E71")
$c.NumberFormat = "@"
$c.Value = "7"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("I71")
$c.NumberFormat = "@"
$c.Value = "7"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("J71")
$c.NumberFormat = "@"
$c.Value = "6"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("E72")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("I72")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("J72")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("N72")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("E73")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("I73")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("J73")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("E74")
$c.NumberFormat = "@"
$c.Value = "6"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("I74")
$c.NumberFormat = "@"
$c.Value = "6"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("J74")
$c.NumberFormat = "@"
$c.Value = "3"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("N74")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("E76")
$c.NumberFormat = "@"
$c.Value = "9"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("I76")
$c.NumberFormat = "@"
$c.Value = "9"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("J76")
$c.NumberFormat = "@"
$c.Value = "10"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("N76")
$c.NumberFormat = "@"
$c.Value = "0"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $false

$c = $ws.Range("E77")
$c.NumberFormat = "@"
$c.Value = "60"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("I77")
$c.NumberFormat = "@"
$c.Value = "60"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("J77")
$c.NumberFormat = "@"
$c.Value = "59"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

$c = $ws.Range("N77")
$c.NumberFormat = "@"
$c.Value = "3"
$c.Font.Name = "Arial"
$c.Font.Size = 9
$c.Font.Bold = $true

